$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 186; this shifts the existing rows 186-214
# down to 187-215, preserving all of their data/formatting.
$ws.Rows(186).Insert()

# Populate the newly inserted row 186 with the new weekly record
# (same Mercado/Region/Categoria/etc. as the surrounding rows).
$ws.Cells.Item(186, 1).Value2 = 3
$ws.Cells.Item(186, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(186, 3).Value2 = "Coquimbo"
$ws.Cells.Item(186, 4).Value2 = 44505
$ws.Cells.Item(186, 5).Value2 = 5
$ws.Cells.Item(186, 6).Value2 = 100112012
$ws.Cells.Item(186, 7).Value2 = "Espinaca"
$ws.Cells.Item(186, 8).Value2 = "Sin especificar"
$ws.Cells.Item(186, 9).Value2 = "Primera"
$ws.Cells.Item(186, 10).Value2 = 250
$ws.Cells.Item(186, 11).Value2 = 2000
$ws.Cells.Item(186, 12).Value2 = 2200
$ws.Cells.Item(186, 13).Value2 = 2104
$ws.Cells.Item(186, 14).Value2 = "$/docena de atados (3 kilos)"
$ws.Cells.Item(186, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(186, 16).Value2 = 701
$ws.Cells.Item(186, 17).Value2 = 3
$ws.Cells.Item(186, 18).Value2 = "Hortaliza"
